# Update code size repo after Functor refactoring.
# Adds a second "ATmega328P" results block (columns J:K) mirroring the
# existing C:D block, refreshes the C:D figures for the first four
# examples with the new compiler's numbers, fills in the previously
# empty J:K figures for every example row, adds a trailing formatted
# pair of cells on the totals row, and drops an explanatory comment on
# the new header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. New header cells for the J:K ("ATmega328P") block, mirroring C:D.
# ---------------------------------------------------------------------

# Merge J1:K1 first, then clone the formatting of the existing C1:D1
# merged header so the new block visually matches the other headers.
$ws.Range("J1:K1").Merge() | Out-Null

$ws.Range("C1").Copy($null) | Out-Null
$ws.Range("J1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D1").Copy($null) | Out-Null
$ws.Range("K1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J1").Value = "ATmega328P"

# Second header row: "Code size" / "Data size" labels, formatted like
# the matching C2/D2 cells.
$ws.Range("J2").Value = "Code size"
$ws.Range("K2").Value = "Data size"
$ws.Range("C2").Copy($null) | Out-Null
$ws.Range("J2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D2").Copy($null) | Out-Null
$ws.Range("K2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J2").Value = "Code size"
$ws.Range("K2").Value = "Data size"

# ---------------------------------------------------------------------
# 2. Refreshed code/data size figures for the first four examples
#    (new compiler run).
# ---------------------------------------------------------------------
$ws.Range("C3").Value = 1174
$ws.Range("D3").Value = 13
$ws.Range("C4").Value = 1014
$ws.Range("D4").Value = 110
$ws.Range("C5").Value = 1332
$ws.Range("D5").Value = 116
$ws.Range("C6").Value = 1338
$ws.Range("D6").Value = 116

# ---------------------------------------------------------------------
# 3. J:K figures for every example (the previous ATmega328P numbers,
#    now carried over verbatim into the new block).
# ---------------------------------------------------------------------
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 1
$ws.Range("J4").Value = 978
$ws.Range("K4").Value = 104
$ws.Range("J5").Value = 1256
$ws.Range("K5").Value = 104
$ws.Range("J6").Value = 1262
$ws.Range("K6").Value = 104
$ws.Range("J7").Value = 380
$ws.Range("K7").Value = 7
$ws.Range("J8").Value = 156
$ws.Range("K8").Value = 0
$ws.Range("J9").Value = 162
$ws.Range("K9").Value = 0
$ws.Range("J10").Value = 232
$ws.Range("K10").Value = 0
$ws.Range("J11").Value = 346
$ws.Range("K11").Value = 0
$ws.Range("J12").Value = 1884
$ws.Range("K12").Value = 169
$ws.Range("J13").Value = 1628
$ws.Range("K13").Value = 158
$ws.Range("J14").Value = 1626
$ws.Range("K14").Value = 158
$ws.Range("J15").Value = 514
$ws.Range("K15").Value = 8
$ws.Range("J16").Value = 704
$ws.Range("K16").Value = 8
$ws.Range("J17").Value = 486
$ws.Range("K17").Value = 8

# ---------------------------------------------------------------------
# 4. Extend the blank "totals" row formatting to the new J:K columns.
# ---------------------------------------------------------------------
$ws.Range("A18").Copy($null) | Out-Null
$ws.Range("J18").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("K18").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------
# 5. Explanatory comment on the new header cell.
# ---------------------------------------------------------------------
$ws.Range("J1").AddComment("With functors.") | Out-Null

# ---------------------------------------------------------------------
# 6. Leave the selection where the author left off.
# ---------------------------------------------------------------------
$ws.Range("M15").Select() | Out-Null
